# Insert a new data row above row 32 (pushing existing rows 32-139 down to
# 33-140) and populate it with a new "Tuna / Segunda" price observation from
# "Provincia de Copiapó", per the commit "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("32").Insert()

$ws.Range("A32").Value = 4
$ws.Range("B32").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C32").Value = "Los Lagos"
$ws.Range("D32").Value = 44526
$ws.Range("E32").Value = 10
$ws.Range("F32").Value = 100112027
$ws.Range("G32").Value = "Melón"
$ws.Range("H32").Value = "Tuna"
$ws.Range("I32").Value = "Segunda"
$ws.Range("J32").Value = 80
$ws.Range("K32").Value = 25000
$ws.Range("L32").Value = 25000
$ws.Range("M32").Value = 25000
$ws.Range("N32").Value = "$/caja 24 unidades"
$ws.Range("O32").Value = "Provincia de Copiapó"
$ws.Range("P32").Value = 1042
$ws.Range("Q32").Value = 24
$ws.Range("R32").Value = "Hortaliza"
